$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole used data block one column to the left: B1:F5 -> A1:E5
# (column A previously held the per-row GENE index, which is dropped since
#  that same data already lives in what becomes the new column E)
$src = $ws.Range("B1:F5")
$dst = $ws.Range("A1:E5")
$dst.Value2 = $src.Value2

# Drop the now-empty column F entirely
$ws.Range("F1:F5").Clear()

# The old column A (rows 2-5) carried the bordered/bold header-like style;
# the new column A data rows should be plain, unformatted values
$ws.Range("A2:A5").ClearFormats()

# Re-apply the bold / centered / bordered header formatting across the new
# header row A1:E1 (matches the style previously used on B1:F1)
$hdr = $ws.Range("A1:E1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1

$ws.Range("A1").Select()
